$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new contact row (row 6): Name in column A, email in column B
$ws.Range("A6").Value = "Mahmudul"
$ws.Range("B6").Value = "mahmudul5809@gmail.com"

# Turn the e-mail address into a mailto hyperlink, matching the other rows
$ws.Hyperlinks.Add($ws.Range("B6"), "mailto:mahmudul5809@gmail.com", [Type]::Missing, [Type]::Missing, "mahmudul5809@gmail.com")

# Keep the new cell's look consistent with a plain (non hyperlink-styled) cell,
# matching the default formatting used for the rest of column A
$ws.Range("B6").Font.Name = $ws.Range("A6").Font.Name
$ws.Range("B6").Font.Underline = $ws.Range("A6").Font.Underline
$ws.Range("B6").Font.Color = $ws.Range("A6").Font.Color

# Move the active selection down to B7, below the newly added row
$ws.Range("B7").Select() | Out-Null
